# New crime data collected - update weekly CompStat figures for the
# 123rd Precinct report (week of 1/30/2023 - 2/5/2023).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -----------------------------------------------
# "Volume 30   Number  4" -> "Volume 30   Number  5"
$ws.Range("A8").Value = "Volume 30   Number  5"

# "Report Covering the Week  1/23/2023  Through  1/29/2023"
#   -> "Report Covering the Week  1/30/2023  Through  2/5/2023"
$ws.Range("C9").Value = "Report Covering the Week  1/30/2023  Through  2/5/2023"

# --- Crime statistics table (rows 16-27) --------------------------------
# Number formats used by the table's numeric styles, so cells that flip
# between the "no data" placeholder (text "0"/"***.*") and an actual
# figure pick up the right style alongside the new value.
$fmtInt = "#,##0"
$fmtPct = "#,##0.0;""-""#,##0.0"
$fmtGeneral = "General"

# Row 16 - Murder (actual figure -> back to "no data" placeholder)
$ws.Range("F16").Value = "0"
$ws.Range("F16").NumberFormat = $fmtGeneral

# Row 17 - Rape
$ws.Range("D17").Value = 2
$ws.Range("D17").NumberFormat = $fmtInt
$ws.Range("E17").Value = -50
$ws.Range("E17").NumberFormat = $fmtPct
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 6
$ws.Range("J17").Value = 6
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 200
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = 50

# Row 18 - Robbery
$ws.Range("D18").Value = 1
$ws.Range("D18").NumberFormat = $fmtInt
$ws.Range("E18").Value = 0
$ws.Range("E18").NumberFormat = $fmtPct
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 1
$ws.Range("G18").NumberFormat = $fmtInt
$ws.Range("H18").Value = 300
$ws.Range("H18").NumberFormat = $fmtPct
$ws.Range("I18").Value = 7
$ws.Range("J18").Value = 1
$ws.Range("J18").NumberFormat = $fmtInt
$ws.Range("K18").Value = 600
$ws.Range("K18").NumberFormat = $fmtPct
$ws.Range("L18").Value = 16.666666666666
$ws.Range("N18").Value = -72

# Row 19 - Fel. Assault
$ws.Range("C19").Value = 7
$ws.Range("E19").Value = 250
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = -5
$ws.Range("I19").Value = 25
$ws.Range("J19").Value = 25
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 56.25
$ws.Range("M19").Value = 56.25
$ws.Range("N19").Value = 108.333333333333

# Row 20 - Burglary
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 10
$ws.Range("K20").Value = 11.111111111111
$ws.Range("L20").Value = 150
$ws.Range("M20").Value = 25
$ws.Range("N20").Value = -85.294117647058

# Row 21 - Gr. Larceny (bold/total style row)
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 100
$ws.Range("F21").Value = 34
$ws.Range("G21").Value = 32
$ws.Range("H21").Value = 6.25
$ws.Range("I21").Value = 49
$ws.Range("J21").Value = 41
$ws.Range("K21").Value = 19.512195121951
$ws.Range("L21").Value = 75
$ws.Range("M21").Value = 32.432432432432
$ws.Range("N21").Value = -56.25

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 12
$ws.Range("E24").Value = 140
$ws.Range("F24").Value = 44
$ws.Range("G24").Value = 23
$ws.Range("H24").Value = 91.304347826087
$ws.Range("I24").Value = 57
$ws.Range("J24").Value = 26
$ws.Range("K24").Value = 119.230769230769
$ws.Range("L24").Value = 185
$ws.Range("M24").Value = -6.557377049180

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 16
$ws.Range("J25").Value = 16
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 700
$ws.Range("M25").Value = -23.809523809523

# Row 26 - Other Sex Crimes
$ws.Range("L26").Value = -100
$ws.Range("L26").NumberFormat = $fmtPct

# Row 27 - Shooting Vic.
$ws.Range("D27").Value = 2
$ws.Range("D27").NumberFormat = $fmtInt
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = $fmtPct
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = -66.666666666666
